$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion message text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.99 = 27422.66 pesos`n✅ 27422.66 pesos = 6.96 = 972.31 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 143
$ws2.Range("O10").Value = 3921.44
$ws2.Range("N12").Value = 3940
$ws2.Range("O12").Value = 139.699
